{"js": "const replacements = [\n  [\"2024-08-01 Thursday\", \"2024-08-02 Friday\"],\n  [\"23\u00f77=3, 2\", \"76\u00f77=10, 6\"],\n  [\"40\u00f78=5, 0\", \"48\u00f73=16, 0\"],\n  [\"97\u00f78=12, 1\", \"24\u00f74=6, 0\"],\n  [\"82\u00f78=10, 2\", \"87\u00f75=17, 2\"],\n  [\"98\u00f72=49, 0\", \"76\u00f76=12, 4\"],\n  [\"91\u00f72=45, 1\", \"19\u00f75=3, 4\"],\n  [\"18\u00f72=9, 0\", \"50\u00f76=8, 2\"],\n  [\"69\u00f77=9, 6\", \"79\u00f78=9, 7\"],\n  [\"85\u00f77=12, 1\", \"23\u00f76=3, 5\"],\n  [\"67\u00f72=33, 1\", \"76\u00f74=19, 0\"],\n  [\"37\u00f75=7, 2\", \"85\u00f72=42, 1\"],\n  [\"92\u00f77=13, 1\", \"75\u00f76=12, 3\"],\n  [\"32\u00f77=4, 4\", \"43\u00f79=4, 7\"],\n  [\"86\u00f78=10, 6\", \"76\u00f76=12, 4\"],\n  [\"17\u00f79=1, 8\", \"67\u00f76=11, 1\"],\n  [\"83\u00f79=9, 2\", \"25\u00f77=3, 4\"],\n  [\"12\u00f74=3, 0\", \"69\u00f72=34, 1\"],\n  [\"43\u00f74=10, 3\", \"91\u00f76=15, 1\"],\n  [\"52\u00f73=17, 1\", \"41\u00f78=5, 1\"],\n  [\"94\u00f78=11, 6\", \"56\u00f72=28, 0\"],\n  [\"84\u00f75=16, 4\", \"34\u00f75=6, 4\"],\n  [\"81\u00f77=11, 4\", \"88\u00f79=9, 7\"],\n  [\"13\u00f78=1, 5\", \"26\u00f75=5, 1\"],\n  [\"75\u00f75=15, 0\", \"96\u00f79=10, 6\"],\n  [\"50\u00f72=25, 0\", \"89\u00f76=14, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2024-08-01 Thursday\", \"2024-08-02 Friday\")\n    ,@(\"23\u00f77=3, 2\", \"76\u00f77=10, 6\")\n    ,@(\"40\u00f78=5, 0\", \"48\u00f73=16, 0\")\n    ,@(\"97\u00f78=12, 1\", \"24\u00f74=6, 0\")\n    ,@(\"82\u00f78=10, 2\", \"87\u00f75=17, 2\")\n    ,@(\"98\u00f72=49, 0\", \"76\u00f76=12, 4\")\n    ,@(\"91\u00f72=45, 1\", \"19\u00f75=3, 4\")\n    ,@(\"18\u00f72=9, 0\", \"50\u00f76=8, 2\")\n    ,@(\"69\u00f77=9, 6\", \"79\u00f78=9, 7\")\n    ,@(\"85\u00f77=12, 1\", \"23\u00f76=3, 5\")\n    ,@(\"67\u00f72=33, 1\", \"76\u00f74=19, 0\")\n    ,@(\"37\u00f75=7, 2\", \"85\u00f72=42, 1\")\n    ,@(\"92\u00f77=13, 1\", \"75\u00f76=12, 3\")\n    ,@(\"32\u00f77=4, 4\", \"43\u00f79=4, 7\")\n    ,@(\"86\u00f78=10, 6\", \"76\u00f76=12, 4\")\n    ,@(\"17\u00f79=1, 8\", \"67\u00f76=11, 1\")\n    ,@(\"83\u00f79=9, 2\", \"25\u00f77=3, 4\")\n    ,@(\"12\u00f74=3, 0\", \"69\u00f72=34, 1\")\n    ,@(\"43\u00f74=10, 3\", \"91\u00f76=15, 1\")\n    ,@(\"52\u00f73=17, 1\", \"41\u00f78=5, 1\")\n    ,@(\"94\u00f78=11, 6\", \"56\u00f72=28, 0\")\n    ,@(\"84\u00f75=16, 4\", \"34\u00f75=6, 4\")\n    ,@(\"81\u00f77=11, 4\", \"88\u00f79=9, 7\")\n    ,@(\"13\u00f78=1, 5\", \"26\u00f75=5, 1\")\n    ,@(\"75\u00f75=15, 0\", \"96\u00f79=10, 6\")\n    ,@(\"50\u00f72=25, 0\", \"89\u00f76=14, 5\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $true, $false, $false, $false, $false, $true, 1, $true, [ref]$newText, 2) | Out-Null\n}\n"}
